$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Sobre mim" backlog entry (the FAZENDO row for "Colocar
# conteúdo na Sobre mim") entirely, shifting all subsequent rows up by
# one and dropping the now-unused shared strings for that row's text.
$ws.Rows.Item(11).Delete()

# Move the active selection to where it ended up after the edit.
$ws.Range("I20").Select()
